$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hydro_reservoir")
$ws.Range("I1:L1").FormulaArray = "=gen_technology!B1:E1"
for ($c = 9; $c -le 12; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    Write-Host $c $cell.Address() "formula=" $cell.Formula "text=" $cell.Text
}
